$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting Type/Index/Original/Translation
# (and the data below them) one column to the right.
$ws.Columns("B:B").Insert()

# Populate the new "Variable" column.
$ws.Range("B1").Value = "Variable"
$ws.Range("B2").Value = "e1"

# Set the widths to match the final layout (values chosen so the stored
# OOXML width ends up as close as possible to the target column widths).
$ws.Columns("B:B").ColumnWidth = 18.6666666666667
$ws.Columns("C:C").ColumnWidth = 17.3333333333333
$ws.Columns("F:F").ColumnWidth = 39

$ws.Range("B3").Select()
